$d = $word.ActiveDocument

$old = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$new = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$r = $d.Content
$r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($r.Find.Found) {
    $r.Delete()
    $r.InsertAfter($new)
}
